$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.842.16"
$ws.Range("E2").Value = "  +0.36%  "

$ws.Range("D3").Value = "1.642.09"
$ws.Range("E3").Value = "  +0.19%  "

$ws.Range("E4").Value = "  -0.71%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.90%  "

$ws.Range("E6").Value = "  +1.55%  "

$ws.Range("E7").Value = "  -0.77%  "

$ws.Range("E8").Value = "  +1.68%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0621"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.14%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.81"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.60%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0845"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.10%  "

$ws.Range("D12").Value = "1.872.02"
$ws.Range("E12").Value = "  +0.17%  "

$ws.Range("D13").Value = "1.647.28"
$ws.Range("E13").Value = "  +0.48%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.13"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.27%  "

$ws.Range("E15").Value = "  +1.20%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.69%  "

$ws.Range("D17").Value = "26.858.55"
$ws.Range("E17").Value = "  +0.46%  "

$ws.Range("E18").Value = "  +1.26%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "217.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.05%  "

$ws.Range("E20").Value = "  -0.81%  "

$ws.Range("E21").Value = "  +1.39%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.39%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.44"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.77%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.26%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.57"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.98%  "

$ws.Range("E26").Value = "  -1.04%  "

$ws.Range("E27").Value = "  +4.75%  "

$ws.Range("E28").Value = "  +0.98%  "

$ws.Range("E29").Value = "  +2.32%  "

$ws.Range("E30").Value = "  +2.43%  "

$ws.Range("E31").Value = "  -0.37%  "

$ws.Range("E32").Value = "  +1.12%  "

$ws.Range("E33").Value = "  +1.85%  "

$ws.Range("E34").Value = "  +2.25%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.45"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.17%  "

$ws.Range("D36").Value = "1.241.93"
$ws.Range("E36").Value = "  -1.69%  "

$ws.Range("E37").Value = "  +0.38%  "

$ws.Range("E38").Value = "  +3.49%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.833"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.03%  "

$ws.Range("E40").Value = "  -0.82%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.804"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.48%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.37"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.31%  "

$ws.Range("D43").Value = "1.785.07"
$ws.Range("E43").Value = "  +0.23%  "

$ws.Range("E44").Value = "  -3.70%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "60.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.59%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.59"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.73%  "

$ws.Range("E47").Value = "  +1.10%  "

$ws.Range("D48").Value = "0.0₆0104"
$ws.Range("E48").Value = "  +1.87%  "

$ws.Range("E49").Value = "  -1.18%  "

$ws.Range("E50").Value = "  +2.09%  "

$ws.Range("E51").Value = "  +1.33%  "
